$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 2.5
$ws.Range("Y2").Value = 17
$ws.Range("AH2").Value = 8
$ws.Range("G8").Value = 2.2
$ws.Range("H8").Value = 3.1
$ws.Range("J8").Value = 2.88
$ws.Range("K8").Value = 2.05
$ws.Range("O8").Value = 1.36
$ws.Range("P8").Value = 3
$ws.Range("X8").Value = 9.5
$ws.Range("AA8").Value = 19
$ws.Range("AH8").Value = 9.5
$ws.Range("AO8").Value = 12
$ws.Range("AP8").Value = 23
$ws.Range("AX8").Value = 21
$ws.Range("N9").Value = 10
$ws.Range("BC10").Value = 151
$ws.Range("BD10").Value = 151
$ws.Range("Q11").Value = 1.82
$ws.Range("R11").Value = 1.92
$ws.Range("G12").Value = 1.57
$ws.Range("I12").Value = 5.25
$ws.Range("X12").Value = 8
$ws.Range("Z12").Value = 12
$ws.Range("AO12").Value = 8
$ws.Range("Q13").Value = 1.84
$ws.Range("R13").Value = 1.89
$ws.Range("O24").Value = 1.18
$ws.Range("P24").Value = 4.5
$ws.Range("Q24").Value = 1.62
$ws.Range("R24").Value = 2.25
$ws.Range("I31").Value = 3.25
$ws.Range("M31").Value = 1.05
$ws.Range("N31").Value = 11
$ws.Range("O31").Value = 1.3
$ws.Range("P31").Value = 3.4
$ws.Range("Q31").Value = 2.03
$ws.Range("R31").Value = 1.83
$ws.Range("S31").Value = 1.37
$ws.Range("T31").Value = 2.75
$ws.Range("AC31").Value = 10
$ws.Range("AG31").Value = 251
$ws.Range("AL31").Value = 26
$ws.Range("AN31").Value = 4.33
$ws.Range("AT31").Value = 2.75
$ws.Range("AW31").Value = 5
$ws.Range("AY31").Value = 26
$ws.Range("AZ31").Value = 51
$ws.Range("S32").Value = 1.47
$ws.Range("S34").Value = 1.47
$ws.Range("G39").Value = 2.75
$ws.Range("J39").Value = 3.4
$ws.Range("M39").Value = 1.13
$ws.Range("N39").Value = 6
$ws.Range("Q39").Value = 2.6
$ws.Range("R39").Value = 1.48
$ws.Range("AA39").Value = 26
$ws.Range("AX39").Value = 17
$ws.Range("AZ39").Value = 51
$ws.Range("G42").Value = 2.35
$ws.Range("I42").Value = 3.2
$ws.Range("O42").Value = 1.44
$ws.Range("P42").Value = 2.63
$ws.Range("Q42").Value = 2.35
$ws.Range("R42").Value = 1.57
$ws.Range("AC42").Value = 7
$ws.Range("AH42").Value = 8
$ws.Range("AQ42").Value = 51
$ws.Range("AS42").Value = 251
$ws.Range("Q43").Value = 2.35
$ws.Range("R43").Value = 1.57
$ws.Range("G45").Value = 1.85
$ws.Range("I45").Value = 4.2
$ws.Range("L45").Value = 4.33
$ws.Range("O45").Value = 1.2
$ws.Range("P45").Value = 4.33
$ws.Range("Q45").Value = 1.7
$ws.Range("R45").Value = 2.1
$ws.Range("AB45").Value = 21
$ws.Range("AD45").Value = 7
$ws.Range("AE45").Value = 13
$ws.Range("AO45").Value = 9.5
$ws.Range("AQ45").Value = 29
$ws.Range("AV45").Value = 51
$ws.Range("AX45").Value = 21
$ws.Range("AY45").Value = 26
$ws.Range("AZ45").Value = 67
$ws.Range("BA45").Value = 81
$ws.Range("O46").Value = 1.17
$ws.Range("P46").Value = 5
$ws.Range("Q46").Value = 1.57
$ws.Range("R46").Value = 2.35
$ws.Range("BD46").Value = 176
